$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date moves from 44902 -> 44505, variety/quality/volume/unit updated
$ws.Range("D2").Value = 44505
$ws.Range("K2").Value = 'Californiana(o)'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 100
$ws.Range("Q2").Value = '$/bandeja 10 kilos'

# Row 3: date moves from 44902 -> 44505, volume/prices/unit updated
$ws.Range("D3").Value = 44505
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("S3").Value = 1500

# Row 4: date moves from 44505 -> 44902, variety/quality/volume/unit updated
$ws.Range("D4").Value = 44902
$ws.Range("K4").Value = 'Golden Nugget'
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 60
$ws.Range("Q4").Value = '$/caja 10 kilos'

# Row 5: date moves from 44505 -> 44902, volume/prices/unit updated
$ws.Range("D5").Value = 44902
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range("Q5").Value = '$/caja 10 kilos'
$ws.Range("S5").Value = 1300
